$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.972.97"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.92"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.44"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.63"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.873.88"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.640.67"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("E15").Value = "  +4.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.92"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.971.19"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.67"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.77"
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.70"
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("E26").Value = "  +1.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.73"
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.34"
$ws.Range("E32").Value = "  +2.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.11"
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.404.19"
$ws.Range("E34").Value = "  -5.19%  "
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("E38").Value = "  +1.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.557"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.915"
$ws.Range("E40").Value = "  -4.61%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.88"
$ws.Range("E43").Value = "  +7.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.33"
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("E45").Value = "  +3.15%  "
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.782.30"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.04"
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("E49").Value = "  +1.33%  "
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.63"
$ws.Range("E51").Value = "  -0.96%  "
